# Applies the weekly data reshuffle for "Alcachofa" (artichoke) rows 2-26
# at "Terminal La Palmera de La Serena" (Hortaliza sheet), keeping rows 3 and 9 unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 16)
$ws.Range("D2").Value = 44855
$ws.Range("H2").Value = 'Española'
$ws.Range("J2").Value = 540
$ws.Range("K2").Value = 9500
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 9750
$ws.Range("N2").Value = '$/caja 30 unidades'
$ws.Range("O2").Value = 'Provincia del Elquí'
$ws.Range("P2").Value = 325
$ws.Range("Q2").Value = 30

# Row 4 (was row 23)
$ws.Range("D4").Value = 44427
$ws.Range("H4").Value = 'Madrigal'
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 12500
$ws.Range("N4").Value = '$/caja 40 unidades'
$ws.Range("O4").Value = 'Provincia de Limarí'
$ws.Range("P4").Value = 312
$ws.Range("Q4").Value = 40

# Row 5 (was row 18)
$ws.Range("D5").Value = 44438
$ws.Range("H5").Value = 'Española'
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 11500
$ws.Range("N5").Value = '$/caja 30 unidades'
$ws.Range("O5").Value = 'Provincia del Elquí'
$ws.Range("P5").Value = 383
$ws.Range("Q5").Value = 30

# Row 6 (was row 12)
$ws.Range("D6").Value = 44701
$ws.Range("H6").Value = 'Española'
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 19000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 19500
$ws.Range("N6").Value = '$/caja 30 unidades'
$ws.Range("O6").Value = 'Provincia del Elquí'
$ws.Range("P6").Value = 650
$ws.Range("Q6").Value = 30

# Row 7 (was row 19)
$ws.Range("D7").Value = 44790
$ws.Range("H7").Value = 'Española'
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 14500
$ws.Range("N7").Value = '$/caja 30 unidades'
$ws.Range("O7").Value = 'Provincia de Limarí'
$ws.Range("P7").Value = 483
$ws.Range("Q7").Value = 30

# Row 8 (was row 20)
$ws.Range("D8").Value = 44790
$ws.Range("H8").Value = 'Madrigal'
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 11500
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 11750
$ws.Range("N8").Value = '$/caja 40 unidades'
$ws.Range("O8").Value = 'Provincia del Elquí'
$ws.Range("P8").Value = 294
$ws.Range("Q8").Value = 40

# Row 10 (was row 13)
$ws.Range("D10").Value = 44420
$ws.Range("H10").Value = 'Madrigal'
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14500
$ws.Range("N10").Value = '$/caja 40 unidades'
$ws.Range("O10").Value = 'Provincia de Limarí'
$ws.Range("P10").Value = 362
$ws.Range("Q10").Value = 40

# Row 11 (was row 14)
$ws.Range("D11").Value = 44420
$ws.Range("H11").Value = 'Madrigal'
$ws.Range("J11").Value = 700
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13500
$ws.Range("N11").Value = '$/caja 40 unidades'
$ws.Range("O11").Value = 'Provincia del Elquí'
$ws.Range("P11").Value = 338
$ws.Range("Q11").Value = 40

# Row 12 (was row 25)
$ws.Range("D12").Value = 44426
$ws.Range("H12").Value = 'Española'
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 11500
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 11750
$ws.Range("N12").Value = '$/caja 30 unidades'
$ws.Range("O12").Value = 'Provincia de Limarí'
$ws.Range("P12").Value = 392
$ws.Range("Q12").Value = 30

# Row 13 (was row 26)
$ws.Range("D13").Value = 44426
$ws.Range("H13").Value = 'Madrigal'
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 12500
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12750
$ws.Range("N13").Value = '$/caja 40 unidades'
$ws.Range("O13").Value = 'Provincia de Limarí'
$ws.Range("P13").Value = 319
$ws.Range("Q13").Value = 40

# Row 14 (was row 21)
$ws.Range("D14").Value = 45079
$ws.Range("H14").Value = 'Madrigal'
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 16500
$ws.Range("L14").Value = 17000
$ws.Range("M14").Value = 16750
$ws.Range("N14").Value = '$/caja 40 unidades'
$ws.Range("O14").Value = 'Provincia del Elquí'
$ws.Range("P14").Value = 419
$ws.Range("Q14").Value = 40

# Row 15 (was row 24)
$ws.Range("D15").Value = 45090
$ws.Range("H15").Value = 'Madrigal'
$ws.Range("J15").Value = 340
$ws.Range("K15").Value = 15500
$ws.Range("L15").Value = 16000
$ws.Range("M15").Value = 15750
$ws.Range("N15").Value = '$/caja 40 unidades'
$ws.Range("O15").Value = 'Provincia del Elquí'
$ws.Range("P15").Value = 394
$ws.Range("Q15").Value = 40

# Row 16 (was row 15)
$ws.Range("D16").Value = 44729
$ws.Range("H16").Value = 'Madrigal'
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 16000
$ws.Range("L16").Value = 17000
$ws.Range("M16").Value = 16500
$ws.Range("N16").Value = '$/caja 40 unidades'
$ws.Range("O16").Value = 'Provincia del Elquí'
$ws.Range("P16").Value = 412
$ws.Range("Q16").Value = 40

# Row 17 (was row 11)
$ws.Range("D17").Value = 44767
$ws.Range("H17").Value = 'Madrigal'
$ws.Range("J17").Value = 600
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("N17").Value = '$/caja 40 unidades'
$ws.Range("O17").Value = 'Provincia de Limarí'
$ws.Range("P17").Value = 362
$ws.Range("Q17").Value = 40

# Row 18 (was row 6)
$ws.Range("D18").Value = 45037
$ws.Range("H18").Value = 'Madrigal'
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 16000
$ws.Range("L18").Value = 17000
$ws.Range("M18").Value = 16500
$ws.Range("N18").Value = '$/caja 40 unidades'
$ws.Range("O18").Value = 'Provincia del Elquí'
$ws.Range("P18").Value = 412
$ws.Range("Q18").Value = 40

# Row 19 (was row 2)
$ws.Range("D19").Value = 44784
$ws.Range("H19").Value = 'Madrigal'
$ws.Range("J19").Value = 520
$ws.Range("K19").Value = 11500
$ws.Range("L19").Value = 12000
$ws.Range("M19").Value = 11750
$ws.Range("N19").Value = '$/caja 40 unidades'
$ws.Range("O19").Value = 'Provincia del Elquí'
$ws.Range("P19").Value = 294
$ws.Range("Q19").Value = 40

# Row 20 (was row 22)
$ws.Range("D20").Value = 44498
$ws.Range("H20").Value = 'Española'
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 8500
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = 8750
$ws.Range("N20").Value = '$/caja 30 unidades'
$ws.Range("O20").Value = 'Provincia de Limarí'
$ws.Range("P20").Value = 292
$ws.Range("Q20").Value = 30

# Row 21 (was row 17)
$ws.Range("D21").Value = 45082
$ws.Range("H21").Value = 'Madrigal'
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17500
$ws.Range("N21").Value = '$/caja 40 unidades'
$ws.Range("O21").Value = 'Provincia del Elquí'
$ws.Range("P21").Value = 438
$ws.Range("Q21").Value = 40

# Row 22 (was row 5)
$ws.Range("D22").Value = 45070
$ws.Range("H22").Value = 'Madrigal'
$ws.Range("J22").Value = 360
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 17500
$ws.Range("N22").Value = '$/caja 40 unidades'
$ws.Range("O22").Value = 'Provincia del Elquí'
$ws.Range("P22").Value = 438
$ws.Range("Q22").Value = 40

# Row 23 (was row 4)
$ws.Range("D23").Value = 44858
$ws.Range("H23").Value = 'Española'
$ws.Range("J23").Value = 500
$ws.Range("K23").Value = 9500
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = 9750
$ws.Range("N23").Value = '$/caja 30 unidades'
$ws.Range("O23").Value = 'Provincia del Elquí'
$ws.Range("P23").Value = 325
$ws.Range("Q23").Value = 30

# Row 24 (was row 10)
$ws.Range("D24").Value = 44839
$ws.Range("H24").Value = 'Española'
$ws.Range("J24").Value = 400
$ws.Range("K24").Value = 12000
$ws.Range("L24").Value = 13000
$ws.Range("M24").Value = 12500
$ws.Range("N24").Value = '$/caja 30 unidades'
$ws.Range("O24").Value = 'Provincia del Elquí'
$ws.Range("P24").Value = 417
$ws.Range("Q24").Value = 30

# Row 25 (was row 7)
$ws.Range("D25").Value = 45069
$ws.Range("H25").Value = 'Madrigal'
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 17000
$ws.Range("L25").Value = 18000
$ws.Range("M25").Value = 17500
$ws.Range("N25").Value = '$/caja 40 unidades'
$ws.Range("O25").Value = 'Provincia del Elquí'
$ws.Range("P25").Value = 438
$ws.Range("Q25").Value = 40

# Row 26 (was row 8)
$ws.Range("D26").Value = 44687
$ws.Range("H26").Value = 'Española'
$ws.Range("J26").Value = 400
$ws.Range("K26").Value = 18000
$ws.Range("L26").Value = 19000
$ws.Range("M26").Value = 18500
$ws.Range("N26").Value = '$/caja 30 unidades'
$ws.Range("O26").Value = 'Provincia de Limarí'
$ws.Range("P26").Value = 617
$ws.Range("Q26").Value = 30
